# Insert a new "Match ID" column at the very left of the sheet, shifting all
# existing columns (A:AC -> B:AD) one place to the right, then populate the
# new column: a header label in the label row and a constant match id (8)
# for every data row (including the hidden totals row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column A; Excel shifts formatting/merges/data right automatically.
$ws.Columns.Item(1).Insert() | Out-Null

# Header label (row 3 holds the column captions used by the data rows).
$ws.Range("A3").Value = "Match ID"

# Constant match id for every player row (4-19) and the hidden summary row (20).
$ws.Range("A4:A19").Value = 8
$ws.Range("A20").Value = 8

# Match the bold "label" styling already used for the other header cells in
# that row (and nowhere else), mirroring column A's look for rows 3-19 only.
$ws.Range("A3:A19").Font.Bold = $true

# Leave the new column selected, matching the saved selection state.
$ws.Range("A3:A19").Select() | Out-Null
